$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly-report rows right above the existing row 324.
# This pushes the former rows 324-328 down to 326-330 (their content/formatting
# travels with them automatically).
$ws.Range("A324:A325").EntireRow.Insert()

# --- New row 324 (new weekly data point, same "caja 36 atados" unit as the
#     old row that is now at 326) ---
$ws.Cells.Item(324, 1).Value = 9
$ws.Cells.Item(324, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(324, 3).Value = "Metropolitana"
$ws.Cells.Item(324, 4).Value = 44628
$ws.Cells.Item(324, 5).Value = 13
$ws.Cells.Item(324, 6).Value = 100112044
$ws.Cells.Item(324, 7).Value = "Perejil"
$ws.Cells.Item(324, 8).Value = "Sin especificar"
$ws.Cells.Item(324, 9).Value = "Primera"
$ws.Cells.Item(324, 10).Value = 25
$ws.Cells.Item(324, 11).Value = 7000
$ws.Cells.Item(324, 12).Value = 7000
$ws.Cells.Item(324, 13).Value = 7000
$ws.Cells.Item(324, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(324, 15).Value = "Región Metropolitana"
$ws.Cells.Item(324, 16).Value = 194
$ws.Cells.Item(324, 17).Value = 36
$ws.Cells.Item(324, 18).Value = "Hortaliza"

# --- New row 325 (new weekly data point, "docena de atados" unit as the
#     old row that is now at 327) ---
$ws.Cells.Item(325, 1).Value = 9
$ws.Cells.Item(325, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(325, 3).Value = "Metropolitana"
$ws.Cells.Item(325, 4).Value = 44628
$ws.Cells.Item(325, 5).Value = 13
$ws.Cells.Item(325, 6).Value = 100112044
$ws.Cells.Item(325, 7).Value = "Perejil"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 61
$ws.Cells.Item(325, 11).Value = 14000
$ws.Cells.Item(325, 12).Value = 15000
$ws.Cells.Item(325, 13).Value = 14508
$ws.Cells.Item(325, 14).Value = "$/docena de atados"
$ws.Cells.Item(325, 15).Value = "Región Metropolitana"
$ws.Cells.Item(325, 16).Value = 4836
$ws.Cells.Item(325, 17).Value = 3
$ws.Cells.Item(325, 18).Value = "Hortaliza"
